$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.872.39"
$ws.Range("E2").Value = "  -1.80%  "

$ws.Range("D3").Value = "2.918.22"
$ws.Range("E3").Value = "  -3.32%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'377.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "

$ws.Range("D6").Value = "'101.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.26%  "

$ws.Range("D7").Value = "'0.532"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.66%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.30%  "

$ws.Range("D10").Value = "'36.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.72%  "

$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("D12").Value = "'0.0831"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.33%  "

$ws.Range("D13").Value = "3.383.43"
$ws.Range("E13").Value = "  -3.11%  "

$ws.Range("D14").Value = "'17.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.11%  "

$ws.Range("D15").Value = "'7.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.67%  "

$ws.Range("D16").Value = "2.924.71"
$ws.Range("E16").Value = "  -2.98%  "

$ws.Range("D17").Value = "'0.965"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.14%  "

$ws.Range("D18").Value = "50.820.82"
$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("D19").Value = "'3.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -10.03%  "

$ws.Range("E20").Value = "  -5.86%  "

$ws.Range("D21").Value = "'12.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.95%  "

$ws.Range("D22").Value = "0.0₃0943"
$ws.Range("E22").Value = "  -2.38%  "

$ws.Range("D23").Value = "'67.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.70%  "

$ws.Range("D24").Value = "'260.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").Value = "'2.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "

$ws.Range("D26").Value = "'8.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.25%  "

$ws.Range("D27").Value = "'7.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.112"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.24%  "

$ws.Range("E30").Value = "  -4.92%  "

$ws.Range("D31").Value = "'25.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.41%  "

$ws.Range("D32").Value = "'9.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.47%  "

$ws.Range("D33").Value = "'50.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.53%  "

$ws.Range("D34").Value = "'33.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.50%  "

$ws.Range("E35").Value = "  -2.15%  "

$ws.Range("D36").Value = "'0.0447"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("E38").Value = "  -6.51%  "

$ws.Range("E39").Value = "  -3.76%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.114"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.56%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "'16.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.13%  "

$ws.Range("D42").Value = "'1.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.22%  "

$ws.Range("D43").Value = "'120.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.68%  "

$ws.Range("D44").Value = "'21.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.92%  "

$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("E46").Value = "  -2.87%  "

$ws.Range("D47").Value = "'0.270"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.92%  "

$ws.Range("D48").Value = "1.995.00"
$ws.Range("E48").Value = "  -3.10%  "

$ws.Range("D49").Value = "'3.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.78%  "

$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'56.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.07%  "
